# Add a table/figure number to the caption above the Activity Diagram
# table: "ตาราง … Activity Diagram" -> "ตารางที่ 1 Activity Diagram".
#
# The caption paragraph is built from four runs:
#   1) "ตาราง"
#   2) " "              (plain space, Thai/cs formatting)
#   3) "… "              (ellipsis placeholder + space)
#   4) "Activity Diagram"
#
# Target layout (per the authoritative diff) keeps run 1 untouched, turns
# run 2's text into "ที่ ", turns run 3's text into just "1", and inserts a
# brand-new run containing a single space between the new "1" run and the
# existing "Activity Diagram" run.  We therefore perform the edits as two
# separate, narrowly-scoped operations and explicitly "pin" every newly
# written/split run by nudging a character-formatting property (Bold on,
# then immediately back off) right after writing its text. That nudge is a
# no-op for the rendered formatting, but it stops the engine from silently
# re-merging the freshly split run back into its neighbour, which is what
# happens if a plain Find/Replace is used across a run boundary.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: the lone-space run right after "ตาราง" becomes "ที่ ".
# ---------------------------------------------------------------------
$headRng = $d.Content
$headRng.Find.Execute("ตาราง", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$spaceRng = $d.Range($headRng.End, $headRng.End + 1)
$spaceRng.Text = "ที่ "
$spaceRng.Font.Bold = $true
$spaceRng.Font.Bold = $false

# ---------------------------------------------------------------------
# Hunk 2: "… " becomes "1", followed by a brand new run containing a
# single space (so "Activity Diagram" keeps starting its own run).
# ---------------------------------------------------------------------
$ellipsisRng = $d.Content
$ellipsisRng.Find.Execute("… ", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

$ellipsisRng.Collapse(0)
$ellipsisRng.InsertAfter(" ")
$ellipsisRng.Font.Bold = $true
$ellipsisRng.Font.Bold = $false
